$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8635009527206421
$ws.Range("B1").Value = 1.816274881362915
$ws.Range("C1").Value = 6.739763736724854
$ws.Range("D1").Value = 1.603439092636108
$ws.Range("E1").Value = 0.9237149357795715
